# Update FAST_holdings model data:
#  - Bump the "as of" date in the confidential disclaimer text (A13)
#  - Refresh the Weight/Percent-Change figures in D2:E10
#
# The sheet ships protected, so it must be unprotected before any cell can
# be written, and is re-protected afterwards to restore that state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-04-28 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value = 0.1047138266875779
$ws.Range("E2").Value = -0.006293402777777679

$ws.Range("D3").Value = 0.1092396539881565
$ws.Range("E3").Value = -0.003885048277438186

$ws.Range("D4").Value = 0.1162475569378272
$ws.Range("E4").Value = -0.0004708097928437516

$ws.Range("D5").Value = 0.1376148516294312
$ws.Range("E5").Value = -0.0009489166534870641

$ws.Range("D6").Value = 0.1323581642715088
$ws.Range("E6").Value = -0.002641824249165814

$ws.Range("D7").Value = 0.1408109841158056
$ws.Range("E7").Value = -0.002100238663484411

$ws.Range("D8").Value = 0.1293399858871195
$ws.Range("E8").Value = -0.003984063745019917

$ws.Range("D9").Value = 0.1296749764825733
$ws.Range("E9").Value = -0.008859351177616714

$ws.Range("D10").Value = 0.9999999999999999
$ws.Range("E10").Value = -0.003578261715174169

$ws.Protect()
